$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (14), shifting N..P to O..Q
$ws.Columns.Item(14).Insert()
$ws.Columns.Item(14).ColumnWidth = 9.166666666666666

# Activate the Repayment Schedule sheet and set its selection
$ws.Activate() | Out-Null
$ws.Range("S6").Select() | Out-Null
